$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("contact")
$ws.Range("D2").Value = "PASS"
